$d = $word.ActiveDocument

# Update the date heading (first paragraph)
$d.Content.Find.Execute("2023-07-31 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-08-01 Tuesday", 2) | Out-Null

# New values for each of the 100 table cells, in row-major reading order
$newValues = @(
    "15+48=",
    "44+37=",
    "9+4=",
    "19+45=",
    "50-32=",
    "50-37=",
    "70-69=",
    "13+68=",
    "80-64=",
    "98-69=",
    "9+24=",
    "25+48=",
    "28+7=",
    "23+38=",
    "70-63=",
    "53+8=",
    "72-56=",
    "27+66=",
    "73-67=",
    "40-18=",
    "71-17=",
    "58+8=",
    "6+29=",
    "47+29=",
    "80-9=",
    "74-68=",
    "90-18=",
    "24+49=",
    "91-89=",
    "8+17=",
    "51-47=",
    "35+36=",
    "77+19=",
    "40-26=",
    "8+14=",
    "13+79=",
    "30-18=",
    "67-9=",
    "8+13=",
    "27+67=",
    "17+26=",
    "16+77=",
    "97-39=",
    "71-22=",
    "79+8=",
    "83-39=",
    "37+54=",
    "70-46=",
    "8+24=",
    "82-55=",
    "26+49=",
    "50-11=",
    "39+14=",
    "6+55=",
    "92-5=",
    "54-48=",
    "47+49=",
    "45+17=",
    "92-17=",
    "28+18=",
    "48+6=",
    "78+16=",
    "66+5=",
    "55-38=",
    "38+34=",
    "82-69=",
    "27+36=",
    "46-29=",
    "19+68=",
    "7+24=",
    "61-12=",
    "49+2=",
    "25+39=",
    "50-33=",
    "84+8=",
    "79+18=",
    "88-29=",
    "64-35=",
    "61-17=",
    "81-19=",
    "46+25=",
    "16+49=",
    "68+23=",
    "82-9=",
    "27+29=",
    "27-18=",
    "87+8=",
    "19+69=",
    "43+9=",
    "9+29=",
    "30-19=",
    "7+64=",
    "71-29=",
    "26-7=",
    "98-49=",
    "49+39=",
    "88+7=",
    "7+29=",
    "29+62=",
    "23-17="
)

$t = $d.Tables.Item(1)
$rows = $t.Rows.Count
$cols = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $t.Cell($r, $c).Range.Text = $newValues[$idx]
        $idx = $idx + 1
    }
}

Write-Host "Updated" $idx "cells"
